$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1961414790996785
$ws.Range("C2").Value = 0.5498392282958199
$ws.Range("J2").Value = 0.01607717041800643
$ws.Range("P2").Value = 0.1527331189710611
$ws.Range("S2").Value = 0.08520900321543408
$ws.Range("B3").Value = 0.00546448087431694
$ws.Range("C3").Value = 0.01912568306010929
$ws.Range("J3").Value = 0.02459016393442623
$ws.Range("P3").Value = 0.7923497267759563
$ws.Range("S3").Value = 0.1584699453551913
$ws.Range("J4").Value = 0.05454545454545454
$ws.Range("P4").Value = 0.7363636363636363
$ws.Range("S4").Value = 0.2090909090909091
$ws.Range("B6").Value = 0.06648936170212766
$ws.Range("D6").Value = 0.02127659574468085
$ws.Range("E6").Value = 0.002659574468085106
$ws.Range("F6").Value = 0.07180851063829788
$ws.Range("J6").Value = 0.2074468085106383
$ws.Range("O6").Value = 0.02659574468085106
$ws.Range("Q6").Value = 0.2180851063829787
$ws.Range("R6").Value = 0.05585106382978723
$ws.Range("S6").Value = 0.3297872340425532
$ws.Range("B7").Value = 0.1148325358851675
$ws.Range("D7").Value = 0.02870813397129187
$ws.Range("E7").Value = 0.007177033492822967
$ws.Range("F7").Value = 0.06698564593301436
$ws.Range("J7").Value = 0.1531100478468899
$ws.Range("O7").Value = 0.02392344497607655
$ws.Range("Q7").Value = 0.1674641148325359
$ws.Range("R7").Value = 0.09330143540669857
$ws.Range("S7").Value = 0.3444976076555024
$ws.Range("B8").Value = 0.1050679851668727
$ws.Range("D8").Value = 0.02719406674907293
$ws.Range("E8").Value = 0.001236093943139679
$ws.Range("F8").Value = 0.05438813349814586
$ws.Range("J8").Value = 0.1199011124845488
$ws.Range("O8").Value = 0.02101359703337454
$ws.Range("Q8").Value = 0.2126081582200247
$ws.Range("R8").Value = 0.04944375772558714
$ws.Range("S8").Value = 0.4091470951792336
$ws.Range("B9").Value = 0.1057934508816121
$ws.Range("D9").Value = 0.02518891687657431
$ws.Range("E9").Value = 0.002518891687657431
$ws.Range("F9").Value = 0.04030226700251889
$ws.Range("J9").Value = 0.1385390428211587
$ws.Range("O9").Value = 0.01763224181360202
$ws.Range("Q9").Value = 0.2040302267002519
$ws.Range("R9").Value = 0.07052896725440806
$ws.Range("S9").Value = 0.3954659949622166
$ws.Range("B10").Value = 0.1220439189189189
$ws.Range("D10").Value = 0.02533783783783784
$ws.Range("E10").Value = 0.001266891891891892
$ws.Range("F10").Value = 0.06883445945945946
$ws.Range("J10").Value = 0.1313344594594595
$ws.Range("O10").Value = 0.01097972972972973
$ws.Range("Q10").Value = 0.2166385135135135
$ws.Range("R10").Value = 0.08065878378378379
$ws.Range("S10").Value = 0.3429054054054054
$ws.Range("F11").Value = 0.001485884101040119
$ws.Range("G11").Value = 0.1337295690936107
$ws.Range("J11").Value = 0.08766716196136701
$ws.Range("K11").Value = 0.1961367013372957
$ws.Range("L11").Value = 0.5542347696879644
$ws.Range("S11").Value = 0.02674591381872214
$ws.Range("G12").Value = 0.7413333333333333
$ws.Range("J12").Value = 0.2026666666666667
$ws.Range("K12").Value = 0.008
$ws.Range("L12").Value = 0.02133333333333333
$ws.Range("S12").Value = 0.02666666666666667
$ws.Range("G13").Value = 0.6559139784946236
$ws.Range("J13").Value = 0.2903225806451613
$ws.Range("S13").Value = 0.05376344086021505
$ws.Range("G14").Value = 0.8
$ws.Range("J14").Value = 0.2
$ws.Range("F15").Value = 0.01518987341772152
$ws.Range("H15").Value = 0.1443037974683544
$ws.Range("I15").Value = 0.05316455696202532
$ws.Range("J15").Value = 0.3518987341772152
$ws.Range("K15").Value = 0.1113924050632911
$ws.Range("M15").Value = 0.01772151898734177
$ws.Range("O15").Value = 0.05569620253164557
$ws.Range("S15").Value = 0.2506329113924051
$ws.Range("F16").Value = 0.01284796573875803
$ws.Range("H16").Value = 0.1970021413276231
$ws.Range("I16").Value = 0.08137044967880086
$ws.Range("J16").Value = 0.3683083511777302
$ws.Range("K16").Value = 0.1284796573875803
$ws.Range("M16").Value = 0.0278372591006424
$ws.Range("N16").Value = 0.004282655246252677
$ws.Range("O16").Value = 0.04710920770877945
$ws.Range("S16").Value = 0.132762312633833
$ws.Range("F17").Value = 0.01084598698481562
$ws.Range("H17").Value = 0.1507592190889371
$ws.Range("I17").Value = 0.09869848156182212
$ws.Range("J17").Value = 0.4370932754880694
$ws.Range("K17").Value = 0.1225596529284165
$ws.Range("M17").Value = 0.01735357917570499
$ws.Range("O17").Value = 0.04663774403470716
$ws.Range("S17").Value = 0.1160520607375271
$ws.Range("F18").Value = 0.0125
$ws.Range("H18").Value = 0.18125
$ws.Range("I18").Value = 0.10625
$ws.Range("J18").Value = 0.35625
$ws.Range("K18").Value = 0.103125
$ws.Range("O18").Value = 0.08437500000000001
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.01369257950530035
$ws.Range("H19").Value = 0.2045053003533569
$ws.Range("I19").Value = 0.09231448763250884
$ws.Range("J19").Value = 0.3564487632508834
$ws.Range("K19").Value = 0.125
$ws.Range("M19").Value = 0.02429328621908127
$ws.Range("N19").Value = 0.001325088339222615
$ws.Range("O19").Value = 0.06846289752650177
$ws.Range("S19").Value = 0.1139575971731449
